$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date that was bumped from
# 2023-09-23 (serial 45192) to 2023-10-03 (serial 45202) for every
# data row (rows 2 through 150).
$ws.Range("C2:C150").Value = "2023-10-03"
